$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3061.524
$ws.Range("I40").Value = 3629.2
$ws.Range("K40").Value = 3629.2
$ws.Range("M40").Value = -3454.2
$ws.Range("H64").Value = 3111.3333
$ws.Range("I64").Value = 3001
$ws.Range("J64").Value = 3142.8572
$ws.Range("K64").Value = 3001
$ws.Range("L64").Value = 3142.8572
$ws.Range("M64").Value = -2753
$ws.Range("N64").Value = -3638.8572
$ws.Range("H67").Value = 3111.3333
$ws.Range("I67").Value = 3001
$ws.Range("J67").Value = 3142.8572
$ws.Range("K67").Value = 3001
$ws.Range("L67").Value = 3142.8572
$ws.Range("M67").Value = -2143
$ws.Range("N67").Value = -4858.8572
$ws.Range("H113").Value = 2999.3333
$ws.Range("I113").Value = 2999.3333
$ws.Range("K113").Value = 2999.3333
$ws.Range("M113").Value = 254.6667000000002
$ws.Range("H141").Value = 2663.818
$ws.Range("I141").Value = 1430.25
$ws.Range("K141").Value = 4290.75
$ws.Range("M141").Value = 889.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1282.909
$ws.Range("I2").Value = 479.8889
$ws.Range("K2").Value = 479.8889
$ws.Range("M2").Value = -366.8889
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H45").Value = 2675.8125
$ws.Range("I45").Value = 2458.0715
$ws.Range("K45").Value = 2458.0715
$ws.Range("M45").Value = -2081.0715
$ws.Range("H95").Value = 27781.4
$ws.Range("J95").Value = 27781.4
$ws.Range("L95").Value = 27781.4
$ws.Range("N95").Value = -33273.4
$ws.Range("H97").Value = 839
$ws.Range("I97").Value = 819.2308
$ws.Range("K97").Value = 819.2308
$ws.Range("M97").Value = -323.2308
$ws.Range("H102").Value = 1057.7778
$ws.Range("I102").Value = 752.5
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 752.5
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = 869.5
$ws.Range("N102").Value = -6744
$ws.Range("H116").Value = 1282.909
$ws.Range("I116").Value = 479.8889
$ws.Range("K116").Value = 479.8889
$ws.Range("M116").Value = 1814.1111
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1282.909
$ws.Range("I3").Value = 479.8889
$ws.Range("K3").Value = 479.8889
$ws.Range("M3").Value = -365.8889
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H20").Value = 2361.7
$ws.Range("I20").Value = 2401.889
$ws.Range("K20").Value = 2401.889
$ws.Range("M20").Value = -2154.889
$ws.Range("H22").Value = 309.75
$ws.Range("I22").Value = 309.75
$ws.Range("K22").Value = 309.75
$ws.Range("M22").Value = -136.75
$ws.Range("H99").Value = 2345.111
$ws.Range("I99").Value = 2106.5557
$ws.Range("K99").Value = 2106.5557
$ws.Range("M99").Value = -608.5556999999999
$ws.Range("H105").Value = 2340.889
$ws.Range("I105").Value = 2071.0625
$ws.Range("K105").Value = 2071.0625
$ws.Range("M105").Value = -324.0625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 443
$ws.Range("I16").Value = 451.6
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 451.6
$ws.Range("L16").Value = 400
$ws.Range("M16").Value = -164.6
$ws.Range("N16").Value = -974
$ws.Range("H113").Value = 443
$ws.Range("I113").Value = 451.6
$ws.Range("J113").Value = 400
$ws.Range("K113").Value = 451.6
$ws.Range("L113").Value = 400
$ws.Range("M113").Value = 1718.4
$ws.Range("N113").Value = -4740
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 2419.6
$ws.Range("I21").Value = 198.66667
$ws.Range("J21").Value = 5751
$ws.Range("K21").Value = 596.00001
$ws.Range("L21").Value = 17253
$ws.Range("M21").Value = -423.00001
$ws.Range("N21").Value = -17599
$ws.Range("H82").Value = 14499.5
$ws.Range("I82").Value = 9749
$ws.Range("J82").Value = 19250
$ws.Range("K82").Value = 29247
$ws.Range("L82").Value = 57750
$ws.Range("M82").Value = -28841
$ws.Range("N82").Value = -58562
$ws.Range("H85").Value = 14499.5
$ws.Range("I85").Value = 9749
$ws.Range("J85").Value = 19250
$ws.Range("K85").Value = 29247
$ws.Range("L85").Value = 57750
$ws.Range("M85").Value = -27843
$ws.Range("N85").Value = -60558
$ws.Range("H140").Value = 1295.909
$ws.Range("I140").Value = 925.5
$ws.Range("K140").Value = 2776.5
$ws.Range("M140").Value = 2403.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5516.4165
$ws.Range("I80").Value = 3389.8
$ws.Range("J80").Value = 7035.4287
$ws.Range("K80").Value = 3389.8
$ws.Range("L80").Value = 7035.4287
$ws.Range("M80").Value = -2391.8
$ws.Range("N80").Value = -9031.4287
$ws.Range("H83").Value = 5516.4165
$ws.Range("I83").Value = 3389.8
$ws.Range("J83").Value = 7035.4287
$ws.Range("K83").Value = 16949
$ws.Range("L83").Value = 35177.14350000001
$ws.Range("M83").Value = -11957
$ws.Range("N83").Value = -45161.14350000001
$ws.Range("H102").Value = 3772
$ws.Range("I102").Value = 2249.5
$ws.Range("J102").Value = 4207
$ws.Range("K102").Value = 2249.5
$ws.Range("L102").Value = 4207
$ws.Range("M102").Value = -627.5
$ws.Range("N102").Value = -7451
$ws.Range("H122").Value = 3988
$ws.Range("I122").Value = 2657
$ws.Range("J122").Value = 5984.5
$ws.Range("K122").Value = 7971
$ws.Range("L122").Value = 17953.5
$ws.Range("M122").Value = -5521
$ws.Range("N122").Value = -22853.5
$ws.Range("H126").Value = 5002.5
$ws.Range("I126").Value = 5002.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15007.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12537.5
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H141").Value = 24949.75
$ws.Range("J141").Value = 24949.75
$ws.Range("L141").Value = 24949.75
$ws.Range("N141").Value = -35309.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1947.4445
$ws.Range("I7").Value = 1410.6
$ws.Range("K7").Value = 1410.6
$ws.Range("M7").Value = -1298.6
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()
$ws.Range("H50").Value = 9309
$ws.Range("H55").Value = 758.1429000000001
$ws.Range("I55").Value = 463
$ws.Range("K55").Value = 463
$ws.Range("M55").Value = -290
$ws.Range("H126").Value = 1947.4445
$ws.Range("I126").Value = 1410.6
$ws.Range("K126").Value = 4231.799999999999
$ws.Range("M126").Value = -1761.799999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 42663
$ws.Range("J123").Value = 42663
$ws.Range("L123").Value = 42663
$ws.Range("N123").Value = -52463
$ws.Range("H125").Value = 92766.336
$ws.Range("J125").Value = 92766.336
$ws.Range("L125").Value = 92766.336
$ws.Range("N125").Value = -102606.336
$ws.Range("H132").Value = 3477.2812
$ws.Range("I132").Value = 2647.1428
$ws.Range("K132").Value = 7941.428400000001
$ws.Range("M132").Value = -5411.428400000001
$ws.Range("H136").Value = 1365.6666
$ws.Range("I136").Value = 934.05884
$ws.Range("K136").Value = 2802.17652
$ws.Range("M136").Value = -252.17652
$ws.Range("H140").Value = 31248.5
$ws.Range("J140").Value = 31248.5
$ws.Range("L140").Value = 31248.5
$ws.Range("N140").Value = -41608.5
